$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 56 (old rows 56-76 shift down to 58-78)
$ws.Range("A56:A57").EntireRow.Insert()

# New row 56 data
$ws.Range("A56").Value = 5
$ws.Range("B56").Value = "Macroferia Regional de Talca"
$ws.Range("C56").Value = "Maule"
$ws.Range("D56").Value = 44809
$ws.Range("D56").Style = $ws.Range("D58").Style
$ws.Range("D56").NumberFormat = $ws.Range("D58").NumberFormat
$ws.Range("E56").Value = 7
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100107
$ws.Range("H56").Value = "Otros"
$ws.Range("I56").Value = 100107002
$ws.Range("J56").Value = "Chirimoya"
$ws.Range("K56").Value = "Cultivar IV Región"
$ws.Range("L56").Value = "Especial"
$ws.Range("M56").Value = 150
$ws.Range("N56").Value = 28000
$ws.Range("O56").Value = 28000
$ws.Range("P56").Value = 28000
$ws.Range("Q56").Value = "$/bandeja 10 kilos"
$ws.Range("R56").Value = "Provincia de Limarí"
$ws.Range("S56").Value = 2800
$ws.Range("T56").Value = 10

# New row 57 data
$ws.Range("A57").Value = 5
$ws.Range("B57").Value = "Macroferia Regional de Talca"
$ws.Range("C57").Value = "Maule"
$ws.Range("D57").Value = 44809
$ws.Range("D57").Style = $ws.Range("D58").Style
$ws.Range("D57").NumberFormat = $ws.Range("D58").NumberFormat
$ws.Range("E57").Value = 7
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100107
$ws.Range("H57").Value = "Otros"
$ws.Range("I57").Value = 100107002
$ws.Range("J57").Value = "Chirimoya"
$ws.Range("K57").Value = "Cultivar IV Región"
$ws.Range("L57").Value = "Primera"
$ws.Range("M57").Value = 100
$ws.Range("N57").Value = 25000
$ws.Range("O57").Value = 25000
$ws.Range("P57").Value = 25000
$ws.Range("Q57").Value = "$/bandeja 10 kilos"
$ws.Range("R57").Value = "Provincia de Limarí"
$ws.Range("S57").Value = 2500
$ws.Range("T57").Value = 10
